# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
#
# This script rewrites the "K" column (column G) values on Sheet1 to reflect
# the newly-regenerated figures.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mapping of row number -> new K value (column G)
$kValues = @{
    2  = 0
    3  = 1
    4  = 1
    5  = 0
    6  = 0
    7  = 0
    8  = 7
    9  = 1
    10 = 0
    11 = 2
    12 = 1
    13 = 1
    16 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
